# Trade #36 closed at 2026-02-17 15:23:19 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.76
$summary.Range("B4").Value = -0.24
$summary.Range("B5").Value = -0.13
$summary.Range("B6").Value = 36
$summary.Range("B7").Value = 11
$summary.Range("B9").Value = 30.56

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.76000000000001
$status.Range("D4").Value = 36
$status.Range("E4").Value = -0.24
$status.Range("F4").Value = -0.24
$status.Range("G4").Value = 30.56

# --- New trade row data, common to "All Trades" and "MarketMaking" sheets ---
$tradeRow = @(36, "2026-02-17", "15:23:13", "MarketMaking", "UP", 0.35, 0.38, "CLOSED", 8.571400000000001, 0.03, 99.76000000000001, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column B ("2026-02-17") looks like a date, so Excel's autodetection would
    # otherwise silently turn it into a date serial number. Force it to be
    # stored as literal text, matching the other text columns, then clear the
    # temporary text format back off the cell so no stray style sticks around.
    $ws.Range("B37").NumberFormat = "@"
    $ws.Cells.Item(37, 2).Value = $tradeRow[1]
    $ws.Range("B37").ClearFormats()

    $ws.Cells.Item(37, 1).Value = $tradeRow[0]
    $ws.Cells.Item(37, 3).Value = $tradeRow[2]
    $ws.Cells.Item(37, 4).Value = $tradeRow[3]
    $ws.Cells.Item(37, 5).Value = $tradeRow[4]
    $ws.Cells.Item(37, 6).Value = $tradeRow[5]
    $ws.Cells.Item(37, 7).Value = $tradeRow[6]
    $ws.Cells.Item(37, 8).Value = $tradeRow[7]
    $ws.Cells.Item(37, 9).Value = $tradeRow[8]
    $ws.Cells.Item(37, 10).Value = $tradeRow[9]
    $ws.Cells.Item(37, 11).Value = $tradeRow[10]
    $ws.Cells.Item(37, 12).Value = $tradeRow[11]
    $ws.Cells.Item(37, 13).Value = $tradeRow[12]
    $ws.Cells.Item(37, 14).Value = $tradeRow[13]
    $ws.Cells.Item(37, 15).Value = $tradeRow[14]
    $ws.Cells.Item(37, 16).Value = $tradeRow[15]
    $ws.Cells.Item(37, 17).Value = $tradeRow[16]
}
